$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the FilesTab query text in cell B4: the query used to return
# `File Type` and `Breed` columns which have been removed.
$newFilesQuery = "MATCH (f:file)-->(parent)`nWITH DISTINCT f, parent`nMATCH (f)-[*]->(c:case)<--(demo:demographic)`n MATCH (s:study)<-[*]-(c)<--(diag:diagnosis)`n MATCH (samp:sample)-->(c) `n WHERE samp.specific_sample_pathology IN [`"Not Applicable`"]  `nWITH DISTINCT f, parent, c, demo, diag, s`nRETURN coalesce(f.file_name, '') AS ``File Name``, `n        coalesce(labels(parent)[0], '') AS ``Association``,`n        coalesce(f.file_description, '') AS ``Description``,`n        coalesce(f.file_format, '') AS ``Format``,`n        coalesce(f.file_size, '') AS ``Size``,`n        coalesce(c.case_id, '') AS ``Case ID``, `n        coalesce(diag.disease_term,'') AS Diagnosis , `n        coalesce(s.clinical_study_designation,'') AS ``Study Code``"

$ws.Range("B4").Value = $newFilesQuery

# Move the view so row 4 is scrolled into position and select B4, matching
# the state Excel saved after the edit was made.
$ws.Range("B4").Select()
$excel.ActiveWindow.ScrollRow = 4
$excel.ActiveWindow.ScrollColumn = 1
